$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the VERANICE row before the 004455356 row (row 65), pushing rows down.
$ws.Rows.Item(65).Insert()
$cA = $ws.Cells.Item(65, 1)
$cA.Value = "'005009947"
$cA.ClearFormats()
$ws.Cells.Item(65, 2).Value = "VERANICE"
$ws.Cells.Item(65, 3).Value = 501.94

# Insert the RENATO row before the 004575621 row. That row was originally at
# row 106, but the insert above shifted it down by one to row 107.
$ws.Rows.Item(107).Insert()
$cA2 = $ws.Cells.Item(107, 1)
$cA2.Value = "'004862672"
$cA2.ClearFormats()
$ws.Cells.Item(107, 2).Value = "RENATO"
$ws.Cells.Item(107, 3).Value = 286.36

# Delete the old RENATO 0.02 row. That row was originally at row 292; the two
# inserts above shifted it down by two to row 294.
$ws.Rows.Item(294).Delete()
